$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("C17").Value = 20815.4167831292
$ws.Range("D17").Value = 3187.112024907716
$ws.Range("E17").Value = 345.1140774549249
$ws.Range("F17").Value = 602086.7667823683

# Row 18
$ws.Range("C18").Value = 20815.4167831292
$ws.Range("D18").Value = 7930.244745634812
$ws.Range("E18").Value = 258.2236536985432
$ws.Range("F18").Value = 387237.0420339032

# Row 19
$ws.Range("C19").Value = 20815.4167831292
$ws.Range("D19").Value = 11113.1656350234
$ws.Range("E19").Value = 250.8808009867362
$ws.Range("F19").Value = 331962.8994059802

# Row 20
$ws.Range("C20").Value = 20815.4167831292
$ws.Range("D20").Value = 27639.69999999999
$ws.Range("E20").Value = 269.2379327662534
$ws.Range("F20").Value = 1094086.860592877

# Row 21
$ws.Range("C21").Value = 29031.0083701363
$ws.Range("D21").Value = 12134.7214666349
$ws.Range("E21").Value = 419.7913200834381
$ws.Range("F21").Value = 192827.7297911076
